$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2/D3 were mistakenly holding the zh_CN text (duplicated from column A).
# Fill in the correct ko_KR (Korean) translations that were missing.
$ws.Range("D2").Value = "공격할 때 이상한 습관이 있다고 듀나에게 지적을 받아버린 굼은, 자신의 습관을 고치고자 노력했지만 뜻대로 잘되지 않았다.`n"
$ws.Range("D3").Value = "그녀의 습관은 전부 과거의 상처가 남긴 것들이었다. 이 습관들을 바꾸기 위해선, 어쩌면 조금 더 시간이 필요할지도 모르겠다.`n"
